$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.659.39"
$ws.Range("E2").Value = "  +0.13%  "
$ws.Range("D3").Value = "1.592.66"
$ws.Range("E3").Value = "  +0.24%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.18%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.514"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.27%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("E8").Value = "  -0.22%  "
$ws.Range("E9").Value = "  -1.69%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.44"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.79%  "
$ws.Range("D12").Value = "1.815.54"
$ws.Range("E12").Value = "  +0.29%  "
$ws.Range("D13").Value = "1.654.38"
$ws.Range("E13").Value = "  +4.37%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.02"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.48%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.520"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.38%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.43"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.47%  "
$ws.Range("D17").Value = "26.611.90"
$ws.Range("E17").Value = "  -0.04%  "
$ws.Range("D18").Value = "0.0₃0727"
$ws.Range("E18").Value = "  +0.00%  "
$ws.Range("E19").Value = "  +0.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "207.31"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.06%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.78"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.69%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.23"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.26%  "
$ws.Range("E23").Value = "  -1.53%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.83"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.69%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.55"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.05%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.17"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.47%  "
$ws.Range("E28").Value = "  +0.09%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.21"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.29%  "
$ws.Range("E30").Value = "  -0.58%  "
$ws.Range("E31").Value = "  -0.08%  "
$ws.Range("E32").Value = "  -0.99%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.665"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.57%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.93"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.28%  "
$ws.Range("D35").Value = "1.278.42"
$ws.Range("E35").Value = "  -3.82%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.45"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.10%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.49"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.35%  "
$ws.Range("E38").Value = "  -0.43%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.837"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.36%  "
$ws.Range("E40").Value = "  +0.04%  "
$ws.Range("E41").Value = "  +0.04%  "
$ws.Range("E42").Value = "  +1.30%  "
$ws.Range("E43").Value = "  -0.32%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.17"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.57%  "
$ws.Range("D45").Value = "1.728.36"
$ws.Range("E45").Value = "  +0.32%  "
$ws.Range("E46").Value = "  +8.70%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "89.84"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.60"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.79%  "
$ws.Range("E49").Value = "  +2.85%  "
$ws.Range("E50").Value = "  -0.16%  "
$ws.Range("E51").Value = "  -0.03%  "
